# Update gh-pages output (generated at 456a3b4)
# "F"/"G" column updates on sheet "展览" and the aggregated sheet "全部类型".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 99
$ws1.Range("F4").Value  = 7297
$ws1.Range("G4").Value  = 49
$ws1.Range("F6").Value  = 426
$ws1.Range("F7").Value  = 3801
$ws1.Range("F9").Value  = 538
$ws1.Range("F11").Value = 610
$ws1.Range("F12").Value = 104

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 99
$ws4.Range("F5").Value  = 7297
$ws4.Range("G5").Value  = 49
$ws4.Range("F8").Value  = 426
$ws4.Range("F9").Value  = 3801
$ws4.Range("F11").Value = 538
$ws4.Range("F13").Value = 610
$ws4.Range("F14").Value = 104
